# Commit message: "update utils.py exclude bug"
# The corrections lookup table incorrectly contained an entry mapping the
# misspelling "langus" to replacement "NA" (a leftover/bogus row that the
# exclude-logic in utils.py was tripping on). Remove that row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the row whose "entry" column (A) equals "langus" and whose
# "replacement" column (B) equals "NA", and delete the entire row so that
# every row below shifts up by one (mirrors the target diff exactly).
$targetRow = 0
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $entry = $ws.Cells.Item($r, 1).Value2
    $replacement = $ws.Cells.Item($r, 2).Value2
    if ($entry -eq "langus" -and $replacement -eq "NA") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete()
}

# Restore a sensible view/selection state on the sheet (best effort).
$ws.Range("B129").Select()
$excel.ActiveWindow.ScrollRow = 113
$excel.ActiveWindow.ScrollColumn = 1
